# Update the answer table cells with the new divisor/quotient/remainder text.
# Using Tables(1).Cell(row, col) addressing avoids ambiguity since several
# source strings ("50÷3=16, 2") are repeated but map to different targets.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $r = $cell.Range
    # Trim trailing cell-mark / paragraph-mark characters before setting text
    $r.End = $r.End - 1
    $r.Text = $text
}

Set-CellText 1 1 "63÷2=31, 1"
Set-CellText 1 2 "51÷7=7, 2"
Set-CellText 1 3 "44÷8=5, 4"
Set-CellText 1 4 "90÷8=11, 2"
Set-CellText 1 5 "58÷3=19, 1"

Set-CellText 5 1 "26÷9=2, 8"
Set-CellText 5 2 "73÷3=24, 1"
Set-CellText 5 3 "30÷6=5, 0"
Set-CellText 5 4 "45÷4=11, 1"
Set-CellText 5 5 "37÷9=4, 1"

Set-CellText 9 1 "44÷2=22, 0"
Set-CellText 9 2 "89÷7=12, 5"
Set-CellText 9 3 "50÷7=7, 1"
Set-CellText 9 4 "53÷9=5, 8"
Set-CellText 9 5 "57÷9=6, 3"

Set-CellText 13 1 "33÷8=4, 1"
Set-CellText 13 2 "27÷2=13, 1"
Set-CellText 13 3 "13÷5=2, 3"
Set-CellText 13 4 "94÷5=18, 4"
Set-CellText 13 5 "90÷4=22, 2"

Set-CellText 17 1 "80÷3=26, 2"
Set-CellText 17 2 "44÷9=4, 8"
Set-CellText 17 3 "88÷3=29, 1"
Set-CellText 17 4 "98÷5=19, 3"
Set-CellText 17 5 "38÷6=6, 2"

Write-Output "Done updating table cells."
